$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.1390903333333333
$ws.Range("H2").Value = 0.4172709999999999
$ws.Range("I2").Value = 0.2062392066578425
$ws.Range("J2").Value = 0.2062392066578425
$ws.Range("M2").Value = 4.925988333333333
$ws.Range("N2").Value = 14.777965
$ws.Range("O2").Value = 0.05656988822582037
$ws.Range("P2").Value = 0.05656988822582035
$ws.Range("Q2").Value = 0.6851573592794443
$ws.Range("R2").Value = 6.166416233514998
$ws.Range("S2").Value = 0.01166692886841602
$ws.Range("T2").Value = 0.01166692886841601
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.1390903333333333
$ws.Range("H3").Value = 0.4172709999999999
$ws.Range("I3").Value = 0.2062392066578425
$ws.Range("J3").Value = 0.2062392066578425
$ws.Range("O3").Value = 0.5464678959362861
$ws.Range("P3").Value = 0.5464678959362861
$ws.Range("Q3").Value = 6.618653708772999
$ws.Range("R3").Value = 59.56788337895699
$ws.Range("S3").Value = 0.1127031053218801
$ws.Range("T3").Value = 0.1127031053218801
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.1390903333333333
$ws.Range("H4").Value = 0.4172709999999999
$ws.Range("I4").Value = 0.2062392066578425
$ws.Range("J4").Value = 0.2062392066578425
$ws.Range("M4").Value = 13.80191933333334
$ws.Range("N4").Value = 41.40575800000001
$ws.Range("O4").Value = 0.1585007882996995
$ws.Range("P4").Value = 0.1585007882996994
$ws.Range("Q4").Value = 1.919713560713111
$ws.Range("R4").Value = 17.277422046418
$ws.Range("S4").Value = 0.03268907683357267
$ws.Range("T4").Value = 0.03268907683357265
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.1390903333333333
$ws.Range("H5").Value = 0.4172709999999999
$ws.Range("I5").Value = 0.2062392066578425
$ws.Range("J5").Value = 0.2062392066578425
$ws.Range("M5").Value = 12.18796133333333
$ws.Range("N5").Value = 36.563884
$ws.Range("O5").Value = 0.1399661476381804
$ws.Range("P5").Value = 0.1399661476381803
$ws.Range("Q5").Value = 1.695227604507111
$ws.Range("R5").Value = 15.257048440564
$ws.Range("S5").Value = 0.02886650724785277
$ws.Range("T5").Value = 0.02886650724785277
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.1390903333333333
$ws.Range("H6").Value = 0.4172709999999999
$ws.Range("I6").Value = 0.2062392066578425
$ws.Range("J6").Value = 0.2062392066578425
$ws.Range("M6").Value = 8.576764333333333
$ws.Range("N6").Value = 25.730293
$ws.Range("O6").Value = 0.09849527990001386
$ws.Range("P6").Value = 0.09849527990001385
$ws.Range("Q6").Value = 1.192945010044778
$ws.Range("R6").Value = 10.736505090403
$ws.Range("S6").Value = 0.020313588386121
$ws.Range("T6").Value = 0.02031358838612099
$ws.Range("G7").Value = 0.5183446666666667
$ws.Range("I7").Value = 0.7685867900859908
$ws.Range("J7").Value = 0.7685867900859908
$ws.Range("M7").Value = 4.925988333333333
$ws.Range("N7").Value = 14.777965
$ws.Range("O7").Value = 0.05656988822582037
$ws.Range("P7").Value = 0.05656988822582035
$ws.Range("Q7").Value = 2.553359780645555
$ws.Range("R7").Value = 22.98023802581
$ws.Range("S7").Value = 0.04347886880700656
$ws.Range("T7").Value = 0.04347886880700656
$ws.Range("G8").Value = 0.5183446666666667
$ws.Range("I8").Value = 0.7685867900859908
$ws.Range("J8").Value = 0.7685867900859908
$ws.Range("O8").Value = 0.5464678959362861
$ws.Range("P8").Value = 0.5464678959362861
$ws.Range("Q8").Value = 24.665580764942
$ws.Range("S8").Value = 0.4200080060227154
$ws.Range("T8").Value = 0.4200080060227154
$ws.Range("G9").Value = 0.5183446666666667
$ws.Range("I9").Value = 0.7685867900859908
$ws.Range("J9").Value = 0.7685867900859908
$ws.Range("M9").Value = 13.80191933333334
$ws.Range("N9").Value = 41.40575800000001
$ws.Range("O9").Value = 0.1585007882996995
$ws.Range("P9").Value = 0.1585007882996994
$ws.Range("Q9").Value = 7.15415127619689
$ws.Range("R9").Value = 64.387361485772
$ws.Range("S9").Value = 0.1218216121053652
$ws.Range("T9").Value = 0.1218216121053652
$ws.Range("G10").Value = 0.5183446666666667
$ws.Range("I10").Value = 0.7685867900859908
$ws.Range("J10").Value = 0.7685867900859908
$ws.Range("M10").Value = 12.18796133333333
$ws.Range("N10").Value = 36.563884
$ws.Range("O10").Value = 0.1399661476381804
$ws.Range("P10").Value = 0.1399661476381803
$ws.Range("Q10").Value = 6.317564754672889
$ws.Range("R10").Value = 56.858082792056
$ws.Range("S10").Value = 0.1075761321339309
$ws.Range("T10").Value = 0.1075761321339309
$ws.Range("G11").Value = 0.5183446666666667
$ws.Range("I11").Value = 0.7685867900859908
$ws.Range("J11").Value = 0.7685867900859908
$ws.Range("M11").Value = 8.576764333333333
$ws.Range("N11").Value = 25.730293
$ws.Range("O11").Value = 0.09849527990001386
$ws.Range("P11").Value = 0.09849527990001385
$ws.Range("Q11").Value = 4.445720049440222
$ws.Range("R11").Value = 40.011480444962
$ws.Range("S11").Value = 0.07570217101697287
$ws.Range("T11").Value = 0.07570217101697285
$ws.Range("G12").Value = 0.01697766666666667
$ws.Range("H12").Value = 0.050933
$ws.Range("I12").Value = 0.0251740032561666
$ws.Range("J12").Value = 0.0251740032561666
$ws.Range("M12").Value = 4.925988333333333
$ws.Range("N12").Value = 14.777965
$ws.Range("O12").Value = 0.05656988822582037
$ws.Range("P12").Value = 0.05656988822582035
$ws.Range("Q12").Value = 0.0836317879272222
$ws.Range("R12").Value = 0.7526860913449999
$ws.Range("S12").Value = 0.001424090550397783
$ws.Range("T12").Value = 0.001424090550397782
$ws.Range("G13").Value = 0.01697766666666667
$ws.Range("H13").Value = 0.050933
$ws.Range("I13").Value = 0.0251740032561666
$ws.Range("J13").Value = 0.0251740032561666
$ws.Range("O13").Value = 0.5464678959362861
$ws.Range("P13").Value = 0.5464678959362861
$ws.Range("Q13").Value = 0.8078871748789999
$ws.Range("R13").Value = 7.270984573911
$ws.Range("S13").Value = 0.01375678459169058
$ws.Range("T13").Value = 0.01375678459169058
$ws.Range("G14").Value = 0.01697766666666667
$ws.Range("H14").Value = 0.050933
$ws.Range("I14").Value = 0.0251740032561666
$ws.Range("J14").Value = 0.0251740032561666
$ws.Range("M14").Value = 13.80191933333334
$ws.Range("N14").Value = 41.40575800000001
$ws.Range("O14").Value = 0.1585007882996995
$ws.Range("P14").Value = 0.1585007882996994
$ws.Range("Q14").Value = 0.2343243858015556
$ws.Range("R14").Value = 2.108919472214
$ws.Range("S14").Value = 0.003990099360761607
$ws.Range("T14").Value = 0.003990099360761606
$ws.Range("G15").Value = 0.01697766666666667
$ws.Range("H15").Value = 0.050933
$ws.Range("I15").Value = 0.0251740032561666
$ws.Range("J15").Value = 0.0251740032561666
$ws.Range("M15").Value = 12.18796133333333
$ws.Range("N15").Value = 36.563884
$ws.Range("O15").Value = 0.1399661476381804
$ws.Range("P15").Value = 0.1399661476381803
$ws.Range("Q15").Value = 0.2069231448635555
$ws.Range("R15").Value = 1.862308303772
$ws.Range("S15").Value = 0.003523508256396647
$ws.Range("T15").Value = 0.003523508256396646
$ws.Range("G16").Value = 0.01697766666666667
$ws.Range("H16").Value = 0.050933
$ws.Range("I16").Value = 0.0251740032561666
$ws.Range("J16").Value = 0.0251740032561666
$ws.Range("M16").Value = 8.576764333333333
$ws.Range("N16").Value = 25.730293
$ws.Range("O16").Value = 0.09849527990001386
$ws.Range("P16").Value = 0.09849527990001385
$ws.Range("Q16").Value = 0.1456134459298889
$ws.Range("R16").Value = 1.310521013369
$ws.Range("S16").Value = 0.002479520496919989
$ws.Range("T16").Value = 0.002479520496919989
